$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are formatted as Text so that numeric-looking strings
# (e.g. "399.", "77  ") are preserved exactly as text, matching the source data.
$ws.Range("A2:E67").NumberFormat = "@"

# Set the new cell values for rows 2-67 (category renamed Laptop -> Scanner and
# the product listings reshuffled/updated per the source edit).
$ws.Range("A2").Value = "Scanner"
$ws.Range("B2").Value = "ADS-3100 High-Speed Desktop Scanner | Compact with Scan Speeds of Up to 40ppm"
$ws.Range("C2").Value = "280."
$ws.Range("D2").Value = "86  "
$ws.Range("E2").Value = "4.3 out of 5 stars"
$ws.Range("A3").Value = "Scanner"
$ws.Range("B3").Value = "SmartOffice S30 High Speed A3 Large Format Duplex Document Scanner, with 100-page Auto Document Feeder (ADF). Scan 12” x 17” Size or Legal-Size Document"
$ws.Range("C3").Value = "499."
$ws.Range("D3").Value = "24  "
$ws.Range("E3").Value = "4.2 out of 5 stars"
$ws.Range("A4").Value = "Scanner"
$ws.Range("B4").Value = "3240 A3 Large Flatbed Scanner, 2400 DPI, CIS Sensor, Scan 12`"x 17`" in 4 sec, Frameless, Auto-Scan, Document & Photo & Book Scanner, Design for Library, School and Soho. Supports Windows & Mac"
$ws.Range("C4").Value = "399."
$ws.Range("D4").Value = "77  "
$ws.Range("E4").Value = "3.5 out of 5 stars"
$ws.Range("A5").Value = "Scanner"
$ws.Range("A6").Value = "Scanner"
$ws.Range("A7").Value = "Scanner"
$ws.Range("A8").Value = "Scanner"
$ws.Range("A9").Value = "Scanner"
$ws.Range("A10").Value = "Scanner"
$ws.Range("A11").Value = "Scanner"
$ws.Range("A12").Value = "Scanner"
$ws.Range("A13").Value = "Scanner"
$ws.Range("A14").Value = "Scanner"
$ws.Range("B14").Value = "3120 A3 Flatbed Scanner, 1200 DPI, CIS Sensor, Scan 12`" x 17`" in 8 sec, Frameless, Auto-Scan, Document & Photo & Book Scanner, Design for Library, School and Soho. Supports Windows & Mac"
$ws.Range("C14").Value = "329."
$ws.Range("D14").Value = "70  "
$ws.Range("E14").Value = "3.7 out of 5 stars"
$ws.Range("A15").Value = "Scanner"
$ws.Range("B15").Value = "ADS-4300N Professional Desktop Scanner with Fast Scan Speeds, Duplex, and Networking"
$ws.Range("C15").Value = "349."
$ws.Range("D15").Value = "21  "
$ws.Range("E15").Value = "4.4 out of 5 stars"
$ws.Range("A16").Value = "Scanner"
$ws.Range("B16").Value = "Large Format Flatbed Scanner OS 1180 - A3 / Tabloid/Legal Size scan, Up to 1200 DPI scan Resolution for Blueprints and Document. Design for Library, School and Soho, Support Mac and PC"
$ws.Range("C16").Value = "349."
$ws.Range("D16").Value = "445  "
$ws.Range("E16").Value = "4.1 out of 5 stars"
$ws.Range("A17").Value = "Scanner"
$ws.Range("A18").Value = "Scanner"
$ws.Range("A19").Value = "Scanner"
$ws.Range("A20").Value = "Scanner"
$ws.Range("A21").Value = "Scanner"
$ws.Range("A22").Value = "Scanner"
$ws.Range("A23").Value = "Scanner"
$ws.Range("A24").Value = "Scanner"
$ws.Range("A25").Value = "Scanner"
$ws.Range("A26").Value = "Scanner"
$ws.Range("A27").Value = "Scanner"
$ws.Range("A28").Value = "Scanner"
$ws.Range("A29").Value = "Scanner"
$ws.Range("A30").Value = "Scanner"
$ws.Range("A31").Value = "Scanner"
$ws.Range("A32").Value = "Scanner"
$ws.Range("B32").Value = "PS3140U Duplex Document Scanner, Citrix Ready & Twain Support for PC and Mac, Scan and Save Batch Documents"
$ws.Range("C32").Value = "369."
$ws.Range("D32").Value = "10  "
$ws.Range("E32").Value = "4.4 out of 5 stars"
$ws.Range("A33").Value = "Scanner"
$ws.Range("B33").Value = "OpticPro A320E - A3 CCD Sensor Flatbed Scanner, 12`" x 17`" scan Area with 7.8 Second Speed. Windows,Mac, ICA & Twain Compliant."
$ws.Range("C33").Value = "649."
$ws.Range("D33").Value = "72  "
$ws.Range("E33").Value = "4.2 out of 5 stars"
$ws.Range("A34").Value = "Scanner"
$ws.Range("B34").Value = "Portable Scanner, Photo Scanner for A4 Documents Pictures Pages Texts in 900 Dpi, Flat Scanning, Include 16G SD Card, Wand Document Scanner Uploads Images to Computer Via USB Cable, No Driver"
$ws.Range("C34").Value = "64."
$ws.Range("D34").Value = "1,150  "
$ws.Range("E34").Value = "3.8 out of 5 stars"
$ws.Range("B35").Value = "Pantum M6552NW All in One Laser Printer Scanner Copier Wireless Monochrome Black and White Printer Home Office - Print Copy Scan, Speed Up to 23 ppm, 50-Sheet ADF, 150 Large Paper Capacity"
$ws.Range("C35").Value = "169."
$ws.Range("D35").Value = "2,514  "
$ws.Range("E35").Value = "3.7 out of 5 stars"
$ws.Range("B36").Value = "for Zebra QLN220 & ZQ610 Portable Thermal Printer Carrying Case with Shoulder Strap"
$ws.Range("C36").Value = "35."
$ws.Range("B37").Value = "NIIMBOT Case Compatible D11 Label Maker, Compatible with D110 Label Printers, Holder for Labeler Makers Printer & Labeling Tapes, Box with Mesh Pocket for Labeler Accessories (Case Only)"
$ws.Range("C37").Value = "13."
$ws.Range("B38").Value = "Canon TS9521C All-In-One Wireless Crafting Photo Printer, 12X12 Printing, White, Amazon Dash Replenishment Ready"
$ws.Range("B39").Value = "Canon IP8720 Wireless Printer, AirPrint and Cloud Compatible, Black"
$ws.Range("B40").Value = "Brother FAX-2840 High Speed Mono Laser Fax Machine, Dark/Light Gray - FAX2840"
$ws.Range("B41").Value = "Brother Printer Wireless, Fast Electronic Label (QL810W), Black"
$ws.Range("B42").Value = "HPRT Portable Printer Wireless Bluetooth Connection MT866 Thermal Inkless Printer Mobile Support Support 8.5`" X 11`" US Letter, Compatible with iPhone Android MacBook Windows Laptops (Upgraded Version)"
$ws.Range("B44").Value = "MUNBYN Bluetooth Label Printer, 130B Wireless Thermal Shipping Printer for 4x6 Shipping Packages Small Business Office or Home, Compatible with iPhone Android iPad Windows macOS Chrome Etsy Ebay"
$ws.Range("B45").Value = "Brother Print & Cut MFC-J1800DW Wireless Color All-in-One Inkjet Printer with Automatic Paper Cutter | Includes 4 Month Refresh Subscription Trial(1), Amazon Dash Replenishment Ready"
$ws.Range("B46").Value = "ASprink 4x6 Direct Thermal Shipping Label Printer, Bluetooth, Support Amazon, Etsy, UPS, FedEx, USPS"
$ws.Range("B47").Value = "KYOCERA 1102V22US0 LASER, COPY,PRINT,SCAN,NET,DUP"
$ws.Range("B48").Value = "300DPI Bluetooth Thermal Label Printer w/Auto Recognition & Rohm Printer Head, Wireless Shipping Label Printer for 1.57`" - 4.25`" Width Labels, Support Windows/MacOS/Linux/Chromebook/Android/iOS"
$ws.Range("C48").Value = "68."
$ws.Range("D48").Value = "121  "
$ws.Range("E48").Value = "4.1 out of 5 stars"
$ws.Range("B49").Value = "Xiaomi Instant Photo Printer 1S Set - Desktop Photo Printer, High-Resolution Image Quality, 6-inch/3-inch Photographic Paper and Ribbon Set, Instant Printing from Smartphone or Computer"
$ws.Range("C49").Value = "129."
$ws.Range("B50").Value = "JADENS Thermal Shipping Label Printer - 4x6 Desktop Label Maker for Shipping Packages, Compatible with Mac, Windows, Work with Ebay, Etsy, Amazon, UPS, Shopify, No Toners."
$ws.Range("B51").Value = "Jiose Thermal Label Printer - Shipping Label Printer for Small Business - Desktop Postage Label Printer Compatible with MacOS, Windows"
$ws.Range("B52").Value = "HPRT Photo Printer 4x6,Wi-Fi Wireless Instant Picture Printer for iPhone, Android, Smartphone, Thermal Dye-Sublimation Printer,AR Video Printing,Portable Photo Printers for Home Use"
$ws.Range("B53").Value = "Brother DCP-L2550DW Wireless All-in-One Monochrome Laser Printer, Black - Print Scan Copy - 2400 x 600 dpi, 36 ppm, 128MB Memory, 250-Sheet, 50-Sheet ADF, Automatic Duplex Printing, Ethernet, Tillsiy"
$ws.Range("B54").Value = "Brother Monochrome Laser Multifunction All-in-One Printer, MFC-L5700DW, Flexible Network Connectivity"
$ws.Range("B57").Value = "Brother HL-L2405W Wireless Compact Monochrome Laser Printer with Mobile Printing, Black & White Output | Includes Refresh Subscription Trial(1), Amazon Dash Replenishment Ready"
$ws.Range("B58").Value = "iDPRT 4X6 Shipping Label Printer, Bluetooth Thermal Label Printer for iPhone/Android, USB Printer for Windows/Mac/Chrome, Suitable for Small Business and Shipping Package, Used for Ebay, UPS, USP"
$ws.Range("B59").Value = "KYOCERA ECOSYS M2540dw All-in-One Monochrome Laser Printer (Print/Copy/Scan/Fax), 42 ppm, Up to Fine 1200dpi, Gigabit Ethernet, USB, Wireless & Wi-Fi Direct, Mobile Print, 5 Line LCD w/Hard Key Panel"
$ws.Range("B62").Value = "POLONO Thermal Label Printer Shipping Label Printer for Shipping Packages, 4x6 Label Printer, Thermal Label Maker, Compatible with Multiple Platforms, Support Multiple Systems"
$ws.Range("B63").Value = "HP Color Laserjet Pro M454dw Single-Function Wireless Laser Printer, White - Print only - 2.7`" Touchscreen, 28 ppm, 600 x 600 dpi, Auto Duplex Printing, 512MB RAM, Ethernet"
$ws.Range("B64").Value = "HP Laserjet Pro 4001 dn Single-Function Monochrome Laser Printer, White - Print only - Mobile Printing, 42 ppm, 1200 x 1200 dpi, Auto 2-Sided Printing, 8.5 x 14, Ethernet, Hi-Speed USB"
$ws.Range("B65").Value = "Xerox C230/DNI Color Printer, Laser, Wireless"
$ws.Range("B66").Value = "BIXOLON SOHO Series Slim Label Maker Printer - Compact Thermal Printer for Shipping Labels - Compatible with Shopify, Ebay, UPS, USPS, FedEx, Amazon & Etsy - 4x6 Label Printers - 6 IPS (XF3-40)"
$ws.Range("C66").Value = "119."
$ws.Range("D66").Value = "30  "
$ws.Range("E66").Value = "3.9 out of 5 stars"
$ws.Range("B67").Value = "BISOFICE Portable 80mm Thermal Label Printer BT Label Maker Sticker Machine with Rechargeable Battery Compatible with iOS Android Computer for Supermarket Clothing Jewelry Retail Store Labeling"
$ws.Range("C67").Value = "39."
$ws.Range("D67").Value = "2  "
$ws.Range("E67").Value = "3.4 out of 5 stars"

# Clear cells that no longer have data for these rows.
$ws.Range("D36").ClearContents()
$ws.Range("E36").ClearContents()
$ws.Range("D37").ClearContents()
$ws.Range("E37").ClearContents()
$ws.Range("D49").ClearContents()
$ws.Range("E49").ClearContents()
